$d = $word.ActiveDocument

$find = "nonparametric location and scale estimator. "
$replace = "nonparametric location and scale estimator in terms of variance and robustness. "

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
